$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Add new header cells: Q1 = "Protein", R1 = "Ontogeny"
# Shared-string table order requires "Ontogeny" (index 27) to be registered
# before "Protein" (index 28), so set R1 first.
$ws.Range("R1").Value = "Ontogeny"
$ws.Range("Q1").Value = "Protein"

# Select Q2 to match the recorded selection in the edited file
$ws.Range("Q2").Select()
